$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "Resolving-Mac" target-cluster rows (old rows 14-17)
$ws.Range("A14:T17").EntireRow.Delete()

# Rewrite rows 2-13 with the refreshed TPM-derived NATMI values
# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Psap"
$ws.Cells.Item(2, 3).Value = "Gpr37l1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 137.0025913333334
$ws.Cells.Item(2, 8).Value = 411.007774
$ws.Cells.Item(2, 9).Value = 0.07043159922291199
$ws.Cells.Item(2, 10).Value = 0.07043159922291199
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.349437
$ws.Cells.Item(2, 14).Value = 1.048311
$ws.Cells.Item(2, 15).Value = 0.3311673788233273
$ws.Cells.Item(2, 16).Value = 0.3311673788233273
$ws.Cells.Item(2, 17).Value = 47.87377450774601
$ws.Cells.Item(2, 18).Value = 430.863970569714
$ws.Cells.Item(2, 19).Value = 0.02332464810098686
$ws.Cells.Item(2, 20).Value = 0.02332464810098686

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Psap"
$ws.Cells.Item(3, 3).Value = "Gpr37l1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 137.0025913333334
$ws.Cells.Item(3, 8).Value = 411.007774
$ws.Cells.Item(3, 9).Value = 0.07043159922291199
$ws.Cells.Item(3, 10).Value = 0.07043159922291199
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.1240676666666667
$ws.Cells.Item(3, 14).Value = 0.372203
$ws.Cells.Item(3, 15).Value = 0.1175810345404931
$ws.Cells.Item(3, 16).Value = 0.1175810345404931
$ws.Cells.Item(3, 17).Value = 16.99759183401356
$ws.Cells.Item(3, 18).Value = 152.978326506122
$ws.Cells.Item(3, 19).Value = 0.008281420300971383
$ws.Cells.Item(3, 20).Value = 0.008281420300971383

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Psap"
$ws.Cells.Item(4, 3).Value = "Gpr37l1"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 137.0025913333334
$ws.Cells.Item(4, 8).Value = 411.007774
$ws.Cells.Item(4, 9).Value = 0.07043159922291199
$ws.Cells.Item(4, 10).Value = 0.07043159922291199
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.5816626666666668
$ws.Cells.Item(4, 14).Value = 1.744988
$ws.Cells.Item(4, 15).Value = 0.5512515866361798
$ws.Cells.Item(4, 16).Value = 0.5512515866361797
$ws.Cells.Item(4, 17).Value = 79.68929261519025
$ws.Cells.Item(4, 18).Value = 717.2036335367121
$ws.Cells.Item(4, 19).Value = 0.03882553082095376
$ws.Cells.Item(4, 20).Value = 0.03882553082095375

# Row 5: FAPs -> ECs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Psap"
$ws.Cells.Item(5, 3).Value = "Gpr37l1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 249.925644
$ws.Cells.Item(5, 8).Value = 749.776932
$ws.Cells.Item(5, 9).Value = 0.1284841594777439
$ws.Cells.Item(5, 10).Value = 0.1284841594777439
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.349437
$ws.Cells.Item(5, 14).Value = 1.048311
$ws.Cells.Item(5, 15).Value = 0.3311673788233273
$ws.Cells.Item(5, 16).Value = 0.3311673788233273
$ws.Cells.Item(5, 17).Value = 87.333267262428
$ws.Cells.Item(5, 18).Value = 785.9994053618519
$ws.Cells.Item(5, 19).Value = 0.04254976231456281
$ws.Cells.Item(5, 20).Value = 0.04254976231456281

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Psap"
$ws.Cells.Item(6, 3).Value = "Gpr37l1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 249.925644
$ws.Cells.Item(6, 8).Value = 749.776932
$ws.Cells.Item(6, 9).Value = 0.1284841594777439
$ws.Cells.Item(6, 10).Value = 0.1284841594777439
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.1240676666666667
$ws.Cells.Item(6, 14).Value = 0.372203
$ws.Cells.Item(6, 15).Value = 0.1175810345404931
$ws.Cells.Item(6, 16).Value = 0.1175810345404931
$ws.Cells.Item(6, 17).Value = 31.007691491244
$ws.Cells.Item(6, 18).Value = 279.069223421196
$ws.Cells.Item(6, 19).Value = 0.01510730039345883
$ws.Cells.Item(6, 20).Value = 0.01510730039345883

# Row 7: FAPs -> MuSCs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Psap"
$ws.Cells.Item(7, 3).Value = "Gpr37l1"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 249.925644
$ws.Cells.Item(7, 8).Value = 749.776932
$ws.Cells.Item(7, 9).Value = 0.1284841594777439
$ws.Cells.Item(7, 10).Value = 0.1284841594777439
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.5816626666666668
$ws.Cells.Item(7, 14).Value = 1.744988
$ws.Cells.Item(7, 15).Value = 0.5512515866361798
$ws.Cells.Item(7, 16).Value = 0.5512515866361797
$ws.Cells.Item(7, 17).Value = 145.372416557424
$ws.Cells.Item(7, 18).Value = 1308.351749016816
$ws.Cells.Item(7, 19).Value = 0.07082709676972229
$ws.Cells.Item(7, 20).Value = 0.07082709676972228

# Row 8: MuSCs -> ECs
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Psap"
$ws.Cells.Item(8, 3).Value = "Gpr37l1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 88.73577866666666
$ws.Cells.Item(8, 8).Value = 266.207336
$ws.Cells.Item(8, 9).Value = 0.04561813567874526
$ws.Cells.Item(8, 10).Value = 0.04561813567874527
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.349437
$ws.Cells.Item(8, 14).Value = 1.048311
$ws.Cells.Item(8, 15).Value = 0.3311673788233273
$ws.Cells.Item(8, 16).Value = 0.3311673788233273
$ws.Cells.Item(8, 17).Value = 31.007564289944
$ws.Cells.Item(8, 18).Value = 279.068078609496
$ws.Cells.Item(8, 19).Value = 0.01510723841953697
$ws.Cells.Item(8, 20).Value = 0.01510723841953697

# Row 9: MuSCs -> FAPs
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Psap"
$ws.Cells.Item(9, 3).Value = "Gpr37l1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 88.73577866666666
$ws.Cells.Item(9, 8).Value = 266.207336
$ws.Cells.Item(9, 9).Value = 0.04561813567874526
$ws.Cells.Item(9, 10).Value = 0.04561813567874527
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.1240676666666667
$ws.Cells.Item(9, 14).Value = 0.372203
$ws.Cells.Item(9, 15).Value = 0.1175810345404931
$ws.Cells.Item(9, 16).Value = 0.1175810345404931
$ws.Cells.Item(9, 17).Value = 11.00924100902311
$ws.Cells.Item(9, 18).Value = 99.083169081208
$ws.Cells.Item(9, 19).Value = 0.005363827586915447
$ws.Cells.Item(9, 20).Value = 0.005363827586915448

# Row 10: MuSCs -> MuSCs
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Psap"
$ws.Cells.Item(10, 3).Value = "Gpr37l1"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 88.73577866666666
$ws.Cells.Item(10, 8).Value = 266.207336
$ws.Cells.Item(10, 9).Value = 0.04561813567874526
$ws.Cells.Item(10, 10).Value = 0.04561813567874527
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.5816626666666668
$ws.Cells.Item(10, 14).Value = 1.744988
$ws.Cells.Item(10, 15).Value = 0.5512515866361798
$ws.Cells.Item(10, 16).Value = 0.5512515866361797
$ws.Cells.Item(10, 17).Value = 51.61428964799645
$ws.Cells.Item(10, 18).Value = 464.528606831968
$ws.Cells.Item(10, 19).Value = 0.02514706967229285
$ws.Cells.Item(10, 20).Value = 0.02514706967229284

# Row 11: Resolving-Mac -> ECs
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Psap"
$ws.Cells.Item(11, 3).Value = "Gpr37l1"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1469.52242
$ws.Cells.Item(11, 8).Value = 4408.56726
$ws.Cells.Item(11, 9).Value = 0.7554661056205989
$ws.Cells.Item(11, 10).Value = 0.7554661056205988
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.349437
$ws.Cells.Item(11, 14).Value = 1.048311
$ws.Cells.Item(11, 15).Value = 0.3311673788233273
$ws.Cells.Item(11, 16).Value = 0.3311673788233273
$ws.Cells.Item(11, 17).Value = 513.50550587754
$ws.Cells.Item(11, 18).Value = 4621.54955289786
$ws.Cells.Item(11, 19).Value = 0.2501857299882406
$ws.Cells.Item(11, 20).Value = 0.2501857299882406

# Row 12: Resolving-Mac -> FAPs
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Psap"
$ws.Cells.Item(12, 3).Value = "Gpr37l1"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1469.52242
$ws.Cells.Item(12, 8).Value = 4408.56726
$ws.Cells.Item(12, 9).Value = 0.7554661056205989
$ws.Cells.Item(12, 10).Value = 0.7554661056205988
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.1240676666666667
$ws.Cells.Item(12, 14).Value = 0.372203
$ws.Cells.Item(12, 15).Value = 0.1175810345404931
$ws.Cells.Item(12, 16).Value = 0.1175810345404931
$ws.Cells.Item(12, 17).Value = 182.3202177637533
$ws.Cells.Item(12, 18).Value = 1640.88195987378
$ws.Cells.Item(12, 19).Value = 0.08882848625914745
$ws.Cells.Item(12, 20).Value = 0.08882848625914744

# Row 13: Resolving-Mac -> MuSCs
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Psap"
$ws.Cells.Item(13, 3).Value = "Gpr37l1"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1469.52242
$ws.Cells.Item(13, 8).Value = 4408.56726
$ws.Cells.Item(13, 9).Value = 0.7554661056205989
$ws.Cells.Item(13, 10).Value = 0.7554661056205988
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.5816626666666668
$ws.Cells.Item(13, 14).Value = 1.744988
$ws.Cells.Item(13, 15).Value = 0.5512515866361798
$ws.Cells.Item(13, 16).Value = 0.5512515866361797
$ws.Cells.Item(13, 17).Value = 854.7663295436535
$ws.Cells.Item(13, 18).Value = 7692.896965892881
$ws.Cells.Item(13, 19).Value = 0.4164518893732109
$ws.Cells.Item(13, 20).Value = 0.4164518893732108
